# doses.xlsx: append 6 new dose rows (two medications x three days) using the
# same source-rows (72:73) formatting as a template, then nudge the view
# state to roughly match what Excel would have recorded after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Copy the formatting of the last existing pair of rows (72:73) down
#        onto the three new pairs of rows (74:75, 76:77, 78:79) so the new
#        cells pick up the same styles (date / integer / boolean / time). ---
$templateRows = $ws.Range("A72:D73")
$templateRows.Copy()
$null = $ws.Range("A74:D75").PasteSpecial(-4122)   # xlPasteFormats
$null = $ws.Range("A76:D77").PasteSpecial(-4122)   # xlPasteFormats
$null = $ws.Range("A78:D79").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- 2. Fill in the new rows' values -----------------------------------
# Row 74: 2024-11-30, medication 8, taken, 10:10
$ws.Cells.Item(74, 1).Value2 = 45626
$ws.Cells.Item(74, 2).Value2 = 8
$ws.Cells.Item(74, 3).Value2 = $true
$ws.Cells.Item(74, 4).Value2 = 0.4236111111111111

# Row 75: 2024-11-30, medication 9, taken, 14:15
$ws.Cells.Item(75, 1).Value2 = 45626
$ws.Cells.Item(75, 2).Value2 = 9
$ws.Cells.Item(75, 3).Value2 = $true
$ws.Cells.Item(75, 4).Value2 = 0.59375

# Row 76: 2024-12-01, medication 8, not taken, 10:10
$ws.Cells.Item(76, 1).Value2 = 45627
$ws.Cells.Item(76, 2).Value2 = 8
$ws.Cells.Item(76, 3).Value2 = $false
$ws.Cells.Item(76, 4).Value2 = 0.4236111111111111

# Row 77: 2024-12-01, medication 9, not taken, 14:15
$ws.Cells.Item(77, 1).Value2 = 45627
$ws.Cells.Item(77, 2).Value2 = 9
$ws.Cells.Item(77, 3).Value2 = $false
$ws.Cells.Item(77, 4).Value2 = 0.59375

# Row 78: 2024-12-02, medication 8, taken, 9:00 (new time format h:mm)
$ws.Cells.Item(78, 1).Value2 = 45628
$ws.Cells.Item(78, 2).Value2 = 8
$ws.Cells.Item(78, 3).Value2 = $true
$ws.Cells.Item(78, 4).Value2 = 0.375
$ws.Cells.Item(78, 4).NumberFormat = "h:mm"

# Row 79: 2024-12-02, medication 9, taken, 14:12 (new time format h:mm)
$ws.Cells.Item(79, 1).Value2 = 45628
$ws.Cells.Item(79, 2).Value2 = 9
$ws.Cells.Item(79, 3).Value2 = $true
$ws.Cells.Item(79, 4).Value2 = 0.59166666666666667
$ws.Cells.Item(79, 4).NumberFormat = "h:mm"

# --- 3. Re-point the view at the newly added rows ------------------------
$null = $ws.Range("I73").Select()

Write-Host "Added rows 74:79 to Sheet1"
